$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 460, shifting existing rows 460:499 down to 461:500.
$ws.Rows.Item(460).Insert()

# Populate the newly inserted row 460 with the new price-report record.
$ws.Cells.Item(460, 1).Value  = 11
$ws.Cells.Item(460, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(460, 3).Value  = "Bíobío"
$ws.Cells.Item(460, 4).Value  = 45106
$ws.Cells.Item(460, 5).Value  = 8
$ws.Cells.Item(460, 6).Value  = 100114001
$ws.Cells.Item(460, 7).Value  = "Papa"
$ws.Cells.Item(460, 8).Value  = "Asterix"
$ws.Cells.Item(460, 9).Value  = "1a (guarda)"
$ws.Cells.Item(460, 10).Value = 3000
$ws.Cells.Item(460, 11).Value = 15000
$ws.Cells.Item(460, 12).Value = 16000
$ws.Cells.Item(460, 13).Value = 15333
$ws.Cells.Item(460, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(460, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(460, 16).Value = 613
$ws.Cells.Item(460, 17).Value = 25
$ws.Cells.Item(460, 18).Value = "Hortaliza"
